$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.223.67'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '1.840.25'
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'232.54"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').Value = "'0.4669"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -3.19%  '
$ws.Range('D8').Value = "'0.2713"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.12%  '
$ws.Range('D9').Value = "'0.06271"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.66%  '
$ws.Range('D10').Value = '1.840.79'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = "'0.07421"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').Value = "'16.09"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.95%  '
$ws.Range('D13').Value = "'4.933"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.00%  '
$ws.Range('D14').Value = "'83.67"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.98%  '
$ws.Range('D15').Value = "'0.6194"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.54%  '
$ws.Range('D16').Value = '30.148.15'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').Value = "'1.000"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = "'225.91"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.93%  '
$ws.Range('D19').Value = "'0.000007271"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.89%  '
$ws.Range('D20').Value = "'12.34"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -5.09%  '
$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').Value = "'1.002"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.081.22'
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('D23').Value = "'4.892"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -5.01%  '
$ws.Range('D24').Value = "'5.848"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.08%  '
$ws.Range('D25').Value = "'9.183"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('D26').Value = "'164.27"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.52%  '
$ws.Range('D27').Value = "'17.75"
$ws.Range('D27').ClearFormats()
$ws.Range('D28').Value = "'1.861"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.31%  '
$ws.Range('E29').Value = '  -1.07%  '
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('D31').Value = "'4.074"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.71%  '
$ws.Range('D32').Value = "'3.807"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.46%  '
$ws.Range('D33').Value = "'0.04818"
$ws.Range('D33').ClearFormats()
$ws.Range('D34').Value = "'1.141"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.22%  '
$ws.Range('D35').Value = "'0.7074"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -4.51%  '
$ws.Range('D36').Value = "'2.704"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('D37').Value = "'0.01869"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.29%  '
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('D39').Value = "'0.8927"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('D40').Value = "'1.918"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -6.35%  '
$ws.Range('E41').Value = '  -1.74%  '
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('D43').Value = "'5.520"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.03%  '
$ws.Range('D44').Value = "'0.4007"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.43%  '
$ws.Range('D45').Value = "'7.029"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.55%  '
$ws.Range('D46').Value = "'0.1192"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.87%  '
$ws.Range('D47').Value = "'59.74"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.69%  '
$ws.Range('D48').Value = "'8.605"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.55%  '
$ws.Range('D49').Value = "'32.88"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.15%  '
$ws.Range('D50').Value = "'0.05513"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.44%  '
$ws.Range('D51').Value = "'1.357"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.62%  '
